# Updates cryptocurrency Price (D) and Volume/1h (E) columns on Sheet1
# with freshly scraped values, per the GitHub Actions commit that refreshed
# the cryptos list. Values that look like plain numbers (single decimal
# point) are written with a leading apostrophe so Excel keeps them as text
# (matching every other cell in these columns) instead of auto-converting
# them to the Number type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '35.112.93'
    'E2' = '  -0.25%  '
    'D3' = '1.904.49'
    'E3' = '  +0.52%  '
    'E4' = '  -0.56%  '
    'D5' = '''253.22'
    'E5' = '  +3.33%  '
    'D6' = '''0.693'
    'E6' = '  +1.35%  '
    'E7' = '  -0.51%  '
    'D8' = '''41.45'
    'E8' = '  +2.89%  '
    'D9' = '''0.357'
    'E9' = '  +3.87%  '
    'D10' = '''52.47'
    'E10' = '  -1.21%  '
    'E11' = '  +4.85%  '
    'E12' = '  -0.80%  '
    'D13' = '''13.18'
    'E13' = '  +5.74%  '
    'D14' = '2.181.57'
    'E14' = '  +0.46%  '
    'D15' = '''0.732'
    'E15' = '  +4.66%  '
    'D16' = '''4.98'
    'E16' = '  +4.76%  '
    'D17' = '1.906.73'
    'E17' = '  +0.40%  '
    'D18' = '35.128.45'
    'E18' = '  -0.28%  '
    'D19' = '''73.57'
    'E19' = '  +2.40%  '
    'D20' = '0.0₃0835'
    'E20' = '  +2.87%  '
    'D21' = '''242.88'
    'E21' = '  +1.38%  '
    'D22' = '''12.91'
    'E22' = '  +3.97%  '
    'D23' = '''5.03'
    'E23' = '  +6.10%  '
    'E24' = '  -0.56%  '
    'E25' = '  +5.08%  '
    'D26' = '''2.29'
    'E26' = '  +0.17%  '
    'D27' = '''167.54'
    'E27' = '  -0.09%  '
    'D28' = '''8.52'
    'E28' = '  +0.86%  '
    'D29' = '''18.50'
    'E29' = '  +2.20%  '
    'E30' = '  +0.14%  '
    'D31' = '4.128.80'
    'E31' = '  -0.33%  '
    'D32' = '''2.04'
    'E32' = '  +11.61%  '
    'E33' = '  +7.91%  '
    'D34' = '''4.32'
    'E34' = '  +4.84%  '
    'D35' = '''1.58'
    'E35' = '  +7.05%  '
    'E36' = '  +3.82%  '
    'E37' = '  -0.49%  '
    'E38' = '  -5.05%  '
    'E39' = '  +0.61%  '
    'D40' = '''103.17'
    'E40' = '  +16.06%  '
    'D41' = '''17.26'
    'E41' = '  +8.50%  '
    'D42' = '''0.0215'
    'E42' = '  +3.84%  '
    'E43' = '  +1.54%  '
    'E44' = '  +2.80%  '
    'D45' = '''2.40'
    'E45' = '  +0.75%  '
    'D46' = '1.309.10'
    'E46' = '  -2.35%  '
    'D47' = '''12.93'
    'E47' = '  +8.41%  '
    'E48' = '  -0.06%  '
    'E49' = '  -1.11%  '
    'D50' = '''6.58'
    'E50' = '  +2.40%  '
    'D51' = '''0.0746'
    'E51' = '  +6.15%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
